# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing N/O/P columns (Disbursement/Paid Date/Outstanding)
# one to the right, then leave that sheet active/selected as the last
# worked-on sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Match the new column's width to its neighbour (column M) before inserting,
# so the freshly-inserted column N keeps a sensible width.
$newColWidth = $ws.Columns("M:M").ColumnWidth

$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $newColWidth

# Make "Repayment schedule" the active sheet/tab with cell R9 selected,
# matching where editing left off.
$ws.Activate() | Out-Null
$ws.Range("R9").Select() | Out-Null
